$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task list")

# Update the "Completed on" date for task #2 (row 3, column F)
$ws.Range("F3").Value = 42036

# Add a new value in the "Progress" column for task #2 (row 3, column G)
$ws.Range("G3").Value = 37302

# Update the active selection on the sheet
[void]$ws.Range("D10").Select()
